$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last status check" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 15.02.2022 09:30"

# Update row 6 (Shell Olomoucka) price data
$ws.Range("B6").Value = 37.9
$ws.Range("C6").Value = 37.5

# D6 and E6 switch from numeric to plain text values with no special
# number formatting - force them to text first so the "+"/date-like
# strings aren't reinterpreted as a number/date, then drop back to the
# default "Normal" style so no numFmt/style index is left on the cell.
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "+0.4"
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2022-02-15 09:33:09"
$ws.Range("E6").Style = "Normal"
